$wb = $excel.ActiveWorkbook

$wsFeatures = $wb.Worksheets.Item("Features")
$wsFeatures.Range("B2").NumberFormat = "@"
$wsFeatures.Range("B2").Value = "0,517"
$wsFeatures.Range("C2").NumberFormat = "@"
$wsFeatures.Range("C2").Value = "0,682"
$wsFeatures.Range("D2").NumberFormat = "@"
$wsFeatures.Range("D2").Value = "0,588"
$wsFeatures.Range("E2").NumberFormat = "@"
$wsFeatures.Range("E2").Value = "0,951"
$wsFeatures.Range("B3").NumberFormat = "@"
$wsFeatures.Range("B3").Value = "0,630"
$wsFeatures.Range("C3").NumberFormat = "@"
$wsFeatures.Range("C3").Value = "0,773"
$wsFeatures.Range("D3").NumberFormat = "@"
$wsFeatures.Range("D3").Value = "0,694"
$wsFeatures.Range("E3").NumberFormat = "@"
$wsFeatures.Range("E3").Value = "0,970"
$wsFeatures.Range("B4").NumberFormat = "@"
$wsFeatures.Range("B4").Value = "0,448"
$wsFeatures.Range("C4").NumberFormat = "@"
$wsFeatures.Range("C4").Value = "0,591"
$wsFeatures.Range("D4").NumberFormat = "@"
$wsFeatures.Range("D4").Value = "0,510"
$wsFeatures.Range("E4").NumberFormat = "@"
$wsFeatures.Range("E4").Value = "0,930"
$wsFeatures.Range("B5").NumberFormat = "@"
$wsFeatures.Range("B5").Value = "0,556"
$wsFeatures.Range("C5").NumberFormat = "@"
$wsFeatures.Range("C5").Value = "0,714"
$wsFeatures.Range("D5").NumberFormat = "@"
$wsFeatures.Range("D5").Value = "0,625"
$wsFeatures.Range("E5").NumberFormat = "@"
$wsFeatures.Range("E5").Value = "0,922"
$wsFeatures.Range("B6").NumberFormat = "@"
$wsFeatures.Range("B6").Value = "0,360"
$wsFeatures.Range("C6").NumberFormat = "@"
$wsFeatures.Range("C6").Value = "0,409"
$wsFeatures.Range("D6").NumberFormat = "@"
$wsFeatures.Range("D6").Value = "0,383"
$wsFeatures.Range("E6").NumberFormat = "@"
$wsFeatures.Range("E6").Value = "0,974"
$wsFeatures.Range("B7").NumberFormat = "@"
$wsFeatures.Range("B7").Value = "0,320"
$wsFeatures.Range("C7").NumberFormat = "@"
$wsFeatures.Range("C7").Value = "0,364"
$wsFeatures.Range("D7").NumberFormat = "@"
$wsFeatures.Range("D7").Value = "0,340"
$wsFeatures.Range("E7").NumberFormat = "@"
$wsFeatures.Range("E7").Value = "1,000"
$wsFeatures.Range("E8").NumberFormat = "@"
$wsFeatures.Range("E8").Value = "0,872"
$wsFeatures.Range("B9").NumberFormat = "@"
$wsFeatures.Range("B9").Value = "0,188"
$wsFeatures.Range("C9").NumberFormat = "@"
$wsFeatures.Range("C9").Value = "0,250"
$wsFeatures.Range("D9").NumberFormat = "@"
$wsFeatures.Range("D9").Value = "0,214"
$wsFeatures.Range("E9").NumberFormat = "@"
$wsFeatures.Range("E9").Value = "0,971"
$wsFeatures.Range("B10").NumberFormat = "@"
$wsFeatures.Range("B10").Value = "0,533"
$wsFeatures.Range("C10").NumberFormat = "@"
$wsFeatures.Range("C10").Value = "0,727"
$wsFeatures.Range("D10").NumberFormat = "@"
$wsFeatures.Range("D10").Value = "0,615"
$wsFeatures.Range("E10").NumberFormat = "@"
$wsFeatures.Range("E10").Value = "0,943"
$wsFeatures.Range("B11").NumberFormat = "@"
$wsFeatures.Range("B11").Value = "0,462"
$wsFeatures.Range("C11").NumberFormat = "@"
$wsFeatures.Range("C11").Value = "0,462"
$wsFeatures.Range("D11").NumberFormat = "@"
$wsFeatures.Range("D11").Value = "0,462"
$wsFeatures.Range("E11").NumberFormat = "@"
$wsFeatures.Range("E11").Value = "0,955"
$wsFeatures.Range("B12").NumberFormat = "@"
$wsFeatures.Range("B12").Value = "0,462"
$wsFeatures.Range("C12").NumberFormat = "@"
$wsFeatures.Range("C12").Value = "0,462"
$wsFeatures.Range("D12").NumberFormat = "@"
$wsFeatures.Range("D12").Value = "0,462"
$wsFeatures.Range("E12").NumberFormat = "@"
$wsFeatures.Range("E12").Value = "0,948"
$wsFeatures.Range("B13").NumberFormat = "@"
$wsFeatures.Range("B13").Value = "0,462"
$wsFeatures.Range("C13").NumberFormat = "@"
$wsFeatures.Range("C13").Value = "0,500"
$wsFeatures.Range("D13").NumberFormat = "@"
$wsFeatures.Range("D13").Value = "0,480"
$wsFeatures.Range("E13").NumberFormat = "@"
$wsFeatures.Range("E13").Value = "0,953"
$wsFeatures.Range("B14").NumberFormat = "@"
$wsFeatures.Range("B14").Value = "0,385"
$wsFeatures.Range("C14").NumberFormat = "@"
$wsFeatures.Range("C14").Value = "0,417"
$wsFeatures.Range("D14").NumberFormat = "@"
$wsFeatures.Range("D14").Value = "0,400"
$wsFeatures.Range("E14").NumberFormat = "@"
$wsFeatures.Range("E14").Value = "0,941"
$wsFeatures.Range("B15").NumberFormat = "@"
$wsFeatures.Range("B15").Value = "0,524"
$wsFeatures.Range("C15").NumberFormat = "@"
$wsFeatures.Range("C15").Value = "0,846"
$wsFeatures.Range("D15").NumberFormat = "@"
$wsFeatures.Range("D15").Value = "0,647"
$wsFeatures.Range("E15").NumberFormat = "@"
$wsFeatures.Range("E15").Value = "0,821"
$wsFeatures.Range("B16").NumberFormat = "@"
$wsFeatures.Range("B16").Value = "0,333"
$wsFeatures.Range("C16").NumberFormat = "@"
$wsFeatures.Range("C16").Value = "0,500"
$wsFeatures.Range("D16").NumberFormat = "@"
$wsFeatures.Range("D16").Value = "0,400"
$wsFeatures.Range("E16").NumberFormat = "@"
$wsFeatures.Range("E16").Value = "0,958"
$wsFeatures.Range("B17").NumberFormat = "@"
$wsFeatures.Range("B17").Value = "0,294"
$wsFeatures.Range("C17").NumberFormat = "@"
$wsFeatures.Range("C17").Value = "0,455"
$wsFeatures.Range("D17").NumberFormat = "@"
$wsFeatures.Range("D17").Value = "0,357"
$wsFeatures.Range("E17").NumberFormat = "@"
$wsFeatures.Range("E17").Value = "0,920"
$wsFeatures.Range("B18").NumberFormat = "@"
$wsFeatures.Range("B18").Value = "0,154"
$wsFeatures.Range("C18").NumberFormat = "@"
$wsFeatures.Range("C18").Value = "0,286"
$wsFeatures.Range("D18").NumberFormat = "@"
$wsFeatures.Range("D18").Value = "0,200"
$wsFeatures.Range("E18").NumberFormat = "@"
$wsFeatures.Range("E18").Value = "0,957"
$wsFeatures.Range("B19").NumberFormat = "@"
$wsFeatures.Range("B19").Value = "0,583"
$wsFeatures.Range("C19").NumberFormat = "@"
$wsFeatures.Range("C19").Value = "0,700"
$wsFeatures.Range("D19").NumberFormat = "@"
$wsFeatures.Range("D19").Value = "0,636"
$wsFeatures.Range("E19").NumberFormat = "@"
$wsFeatures.Range("E19").Value = "0,642"
$wsFeatures.Range("B20").NumberFormat = "@"
$wsFeatures.Range("B20").Value = "0,417"
$wsFeatures.Range("C20").NumberFormat = "@"
$wsFeatures.Range("C20").Value = "0,625"
$wsFeatures.Range("D20").NumberFormat = "@"
$wsFeatures.Range("D20").Value = "0,500"
$wsFeatures.Range("E20").NumberFormat = "@"
$wsFeatures.Range("E20").Value = "0,684"
$wsFeatures.Range("B21").NumberFormat = "@"
$wsFeatures.Range("B21").Value = "0,154"
$wsFeatures.Range("C21").NumberFormat = "@"
$wsFeatures.Range("C21").Value = "0,500"
$wsFeatures.Range("D21").NumberFormat = "@"
$wsFeatures.Range("D21").Value = "0,235"
$wsFeatures.Range("E21").NumberFormat = "@"
$wsFeatures.Range("E21").Value = "0,813"
$wsFeatures.Range("B22").NumberFormat = "@"
$wsFeatures.Range("B22").Value = "0,250"
$wsFeatures.Range("C22").NumberFormat = "@"
$wsFeatures.Range("C22").Value = "0,167"
$wsFeatures.Range("D22").NumberFormat = "@"
$wsFeatures.Range("D22").Value = "0,200"
$wsFeatures.Range("E22").NumberFormat = "@"
$wsFeatures.Range("E22").Value = "0,909"
$wsFeatures.Range("B23").NumberFormat = "@"
$wsFeatures.Range("B23").Value = "0,333"
$wsFeatures.Range("C23").NumberFormat = "@"
$wsFeatures.Range("C23").Value = "1,000"
$wsFeatures.Range("D23").NumberFormat = "@"
$wsFeatures.Range("D23").Value = "0,500"
$wsFeatures.Range("E23").NumberFormat = "@"
$wsFeatures.Range("E23").Value = "0,333"
$wsFeatures.Range("E24").NumberFormat = "@"
$wsFeatures.Range("E24").Value = "0,913"
$wsFeatures.Range("B26").NumberFormat = "@"
$wsFeatures.Range("B26").Value = "0,143"
$wsFeatures.Range("C26").NumberFormat = "@"
$wsFeatures.Range("C26").Value = "0,200"
$wsFeatures.Range("D26").NumberFormat = "@"
$wsFeatures.Range("D26").Value = "0,167"
$wsFeatures.Range("E26").NumberFormat = "@"
$wsFeatures.Range("E26").Value = "1,000"
$wsFeatures.Range("B27").NumberFormat = "@"
$wsFeatures.Range("B27").Value = "1,000"
$wsFeatures.Range("C27").NumberFormat = "@"
$wsFeatures.Range("C27").Value = "0,167"
$wsFeatures.Range("D27").NumberFormat = "@"
$wsFeatures.Range("D27").Value = "0,286"
$wsFeatures.Range("E27").NumberFormat = "@"
$wsFeatures.Range("E27").Value = "0,167"
$wsFeatures.Range("E28").NumberFormat = "@"
$wsFeatures.Range("E28").Value = "1,000"
$wsFeatures.Range("B29").NumberFormat = "@"
$wsFeatures.Range("B29").Value = "0,100"
$wsFeatures.Range("C29").NumberFormat = "@"
$wsFeatures.Range("C29").Value = "0,200"
$wsFeatures.Range("D29").NumberFormat = "@"
$wsFeatures.Range("D29").Value = "0,133"
$wsFeatures.Range("E29").NumberFormat = "@"
$wsFeatures.Range("E29").Value = "1,000"
$wsFeatures.Range("B30").NumberFormat = "@"
$wsFeatures.Range("B30").Value = "0,222"
$wsFeatures.Range("C30").NumberFormat = "@"
$wsFeatures.Range("C30").Value = "0,667"
$wsFeatures.Range("D30").NumberFormat = "@"
$wsFeatures.Range("D30").Value = "0,333"
$wsFeatures.Range("E30").NumberFormat = "@"
$wsFeatures.Range("E30").Value = "0,836"
$wsFeatures.Range("B32").NumberFormat = "@"
$wsFeatures.Range("B32").Value = "0,125"
$wsFeatures.Range("C32").NumberFormat = "@"
$wsFeatures.Range("C32").Value = "0,200"
$wsFeatures.Range("D32").NumberFormat = "@"
$wsFeatures.Range("D32").Value = "0,154"
$wsFeatures.Range("E32").NumberFormat = "@"
$wsFeatures.Range("E32").Value = "1,000"
$wsFeatures.Range("B33").NumberFormat = "@"
$wsFeatures.Range("B33").Value = "1,000"
$wsFeatures.Range("C33").NumberFormat = "@"
$wsFeatures.Range("C33").Value = "1,000"
$wsFeatures.Range("D33").NumberFormat = "@"
$wsFeatures.Range("D33").Value = "1,000"
$wsFeatures.Range("E33").NumberFormat = "@"
$wsFeatures.Range("E33").Value = "1,000"
$wsFeatures.Range("B35").NumberFormat = "@"
$wsFeatures.Range("B35").Value = "1,000"
$wsFeatures.Range("C35").NumberFormat = "@"
$wsFeatures.Range("C35").Value = "0,250"
$wsFeatures.Range("D35").NumberFormat = "@"
$wsFeatures.Range("D35").Value = "0,400"
$wsFeatures.Range("E35").NumberFormat = "@"
$wsFeatures.Range("E35").Value = "0,250"
$wsFeatures.Range("B36").NumberFormat = "@"
$wsFeatures.Range("B36").Value = "0,333"
$wsFeatures.Range("C36").NumberFormat = "@"
$wsFeatures.Range("C36").Value = "0,250"
$wsFeatures.Range("D36").NumberFormat = "@"
$wsFeatures.Range("D36").Value = "0,286"
$wsFeatures.Range("E36").NumberFormat = "@"
$wsFeatures.Range("E36").Value = "1,000"
$wsFeatures.Range("E37").NumberFormat = "@"
$wsFeatures.Range("E37").Value = "1,000"
$wsFeatures.Range("E38").NumberFormat = "@"
$wsFeatures.Range("E38").Value = "1,000"
$wsFeatures.Range("E39").NumberFormat = "@"
$wsFeatures.Range("E39").Value = "1,000"
$wsFeatures.Range("B40").NumberFormat = "@"
$wsFeatures.Range("B40").Value = "0,200"
$wsFeatures.Range("C40").NumberFormat = "@"
$wsFeatures.Range("C40").Value = "0,333"
$wsFeatures.Range("D40").NumberFormat = "@"
$wsFeatures.Range("D40").Value = "0,250"
$wsFeatures.Range("E40").NumberFormat = "@"
$wsFeatures.Range("E40").Value = "1,000"
$wsFeatures.Range("B43").NumberFormat = "@"
$wsFeatures.Range("B43").Value = "0,333"
$wsFeatures.Range("C43").NumberFormat = "@"
$wsFeatures.Range("C43").Value = "0,500"
$wsFeatures.Range("D43").NumberFormat = "@"
$wsFeatures.Range("D43").Value = "0,400"
$wsFeatures.Range("E43").NumberFormat = "@"
$wsFeatures.Range("E43").Value = "1,000"
$wsFeatures.Range("E44").NumberFormat = "@"
$wsFeatures.Range("E44").Value = "1,000"
$wsFeatures.Range("B45").NumberFormat = "@"
$wsFeatures.Range("B45").Value = "0,167"
$wsFeatures.Range("C45").NumberFormat = "@"
$wsFeatures.Range("C45").Value = "0,333"
$wsFeatures.Range("D45").NumberFormat = "@"
$wsFeatures.Range("D45").Value = "0,222"
$wsFeatures.Range("E45").NumberFormat = "@"
$wsFeatures.Range("E45").Value = "1,000"
$wsFeatures.Range("E47").NumberFormat = "@"
$wsFeatures.Range("E47").Value = "1,000"
$wsFeatures.Range("E48").NumberFormat = "@"
$wsFeatures.Range("E48").Value = "1,000"
$wsFeatures.Range("E53").NumberFormat = "@"
$wsFeatures.Range("E53").Value = "0,243"

$wsGlobal = $wb.Worksheets.Item("Global Metrics")
$wsGlobal.Range("B2").NumberFormat = "@"
$wsGlobal.Range("B2").Value = "0,222"
$wsGlobal.Range("C2").NumberFormat = "@"
$wsGlobal.Range("C2").Value = "0,623"
$wsGlobal.Range("D2").NumberFormat = "@"
$wsGlobal.Range("D2").Value = "0,390"
$wsGlobal.Range("E2").NumberFormat = "@"
$wsGlobal.Range("E2").Value = "0,850"
